# Adds a new row (row 10) to Sheet1 containing two date values, as
# described by the commit "Updated the test spreadsheets to include date
# values."
#
#   A10 = "data8"
#   B10 = "row with two dates"
#   C10 = 2012-01-01             (date only)      -> numFmt "MM/DD/YY"
#   D10 = 2012-02-14 02:14:00    (date + time)     -> numFmt "MM/DD/YYYY HH:MM:SS"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "data8"
$ws.Range("B10").Value = "row with two dates"

# Use the raw date serials (days since 1899-12-30, Excel's epoch) so the
# stored value matches exactly, then apply the custom date/time formats
# that show up as new numFmts (165/166) + new cellXfs (2/3) in styles.xml.
$ws.Range("C10").Value = 40909
$ws.Range("C10").NumberFormat = "MM/DD/YY"

$ws.Range("D10").Value = 40953.0930555556
$ws.Range("D10").NumberFormat = "MM/DD/YYYY HH:MM:SS"

# Move the active selection on Sheet1 to A10, matching the new selection.
[void]$ws.Range("A10").Select()
